$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (the "Förändrad" date) for rows 2-7 from 45243 to 45244
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
